$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8624513149261475
$ws.Range("B1").Value = 0.6703342199325562
$ws.Range("C1").Value = 3.491618871688843
$ws.Range("D1").Value = 3.443800210952759
$ws.Range("E1").Value = 0.956577479839325
